# Add a "company_lei" column to the "User input" sheet's Table7, positioned
# right after "company_isin" (so the new column order is:
#   company_name, company_id, company_isin, company_lei, investment_value,
#   engagement_targets, additional_field_1, additional_field_2)
# and populate the two example rows with dummy LEI values, mirroring the
# upstream commit "Added LEI to standard example portfolio".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lo = $ws.ListObjects.Item(1)
$tableName = $lo.Name

# Convert the table back to a plain range first (keeps all existing data &
# formatting intact) so we can freely shift columns around; we rebuild the
# ListObject afterwards.
$lo.Unlist()

# Insert a new blank column in front of the old column D (company_isin's
# neighbour / old investment_value) - this shifts investment_value,
# engagement_targets, additional_field_1 and additional_field_2 one column
# to the right (D->E, E->F, F->G, G->H) together with their data/styles.
$ws.Columns.Item(4).Insert()

# New column header + sample data (LEI = Legal Entity Identifier).
$ws.Range("D1").Value2 = "company_lei"
$ws.Range("D2").Value2 = "JP0000000001"
$ws.Range("D3").Value2 = "UK0000000002"

# Give the new column the same plain width as its neighbour "company_isin"
# column (14 characters, no bestFit autosizing).
$ws.Columns.Item(4).ColumnWidth = 13.17

# Rebuild the table over the new A1:H52 extent with headers, then restore
# the original table/column names (Add() does not read column names back
# from the header cells, so they must be set explicitly).
$newLo = $ws.ListObjects.Add(1, $ws.Range("A1:H52"), [System.Reflection.Missing]::Value, 1)
$newLo.Name = $tableName

$columnNames = @("company_name", "company_id", "company_isin", "company_lei", "investment_value", "engagement_targets", "additional_field_1", "additional_field_2")
for ($i = 1; $i -le $columnNames.Count; $i++) {
    $newLo.ListColumns.Item($i).Name = $columnNames[$i - 1]
}
